# Updated cryptos list - applies price/volume refresh plus a few row swaps
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to text format before writing a numeric-looking
# string so Excel's automatic type detection does not convert it to a
# Number cell (the source data stores these as plain text/inlineStr).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Simple value updates (no row reordering) ---

Set-TextValue $ws.Range("D2") "52.275.34"
$ws.Range("E2").Value  = "  -0.08%  "

Set-TextValue $ws.Range("D3") "2.929.38"
$ws.Range("E3").Value  = "  +0.77%  "

$ws.Range("E4").Value  = "  +0.02%  "

Set-TextValue $ws.Range("D5") "357.94"
$ws.Range("E5").Value  = "  +1.55%  "

Set-TextValue $ws.Range("D6") "110.27"
$ws.Range("E6").Value  = "  -1.47%  "

$ws.Range("E7").Value  = "  +1.59%  "

$ws.Range("E8").Value  = "  +0.00%  "

$ws.Range("E9").Value  = "  +0.38%  "

Set-TextValue $ws.Range("D10") "39.06"
$ws.Range("E10").Value = "  -2.09%  "

$ws.Range("E11").Value = "  +1.45%  "

Set-TextValue $ws.Range("D12") "0.0871"
$ws.Range("E12").Value = "  +0.63%  "

Set-TextValue $ws.Range("D13") "19.62"
$ws.Range("E13").Value = "  -1.20%  "

Set-TextValue $ws.Range("D14") "7.77"
$ws.Range("E14").Value = "  -0.36%  "

Set-TextValue $ws.Range("D15") "3.389.04"
$ws.Range("E15").Value = "  +0.68%  "

Set-TextValue $ws.Range("D16") "2.934.52"
$ws.Range("E16").Value = "  +1.11%  "

Set-TextValue $ws.Range("D17") "0.988"
$ws.Range("E17").Value = "  -1.36%  "

Set-TextValue $ws.Range("D18") "52.270.88"
$ws.Range("E18").Value = "  -0.15%  "

Set-TextValue $ws.Range("D19") "3.56"
$ws.Range("E19").Value = "  +7.73%  "

Set-TextValue $ws.Range("D20") "7.58"
$ws.Range("E20").Value = "  -0.43%  "

Set-TextValue $ws.Range("D21") "13.95"
$ws.Range("E21").Value = "  -1.60%  "

$sub3 = [char]0x2083
$d22Value = "0.0{0}0983" -f $sub3
Set-TextValue $ws.Range("D22") $d22Value
$ws.Range("E22").Value = "  +0.45%  "

Set-TextValue $ws.Range("D23") "70.61"
$ws.Range("E23").Value = "  -0.46%  "

Set-TextValue $ws.Range("D24") "268.70"
$ws.Range("E24").Value = "  -0.41%  "

$ws.Range("E25").Value = "  +1.35%  "

Set-TextValue $ws.Range("D26") "0.184"
$ws.Range("E26").Value = "  +6.81%  "

# --- Row 27 / Row 28 swap (Filecoin <-> EthereumClassic) with updated data ---

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D27") "27.00"
$ws.Range("E27").Value = "  +0.82%  "

$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D28") "7.70"
$ws.Range("E28").Value = "  +15.91%  "

# --- continue simple updates ---

Set-TextValue $ws.Range("D30") "0.106"
$ws.Range("E30").Value = "  +8.53%  "

Set-TextValue $ws.Range("D31") "10.51"
$ws.Range("E31").Value = "  -1.41%  "

Set-TextValue $ws.Range("D32") "37.55"

$ws.Range("E33").Value = "  -2.25%  "

$ws.Range("E34").Value = "  -1.61%  "

Set-TextValue $ws.Range("D35") "52.28"
$ws.Range("E35").Value = "  -2.03%  "

$ws.Range("E36").Value = "  -1.40%  "

# --- Row 38 / 39 / 40 rotation: LidoDAOToken / Celestia / Stacks -> Stacks / LidoDAOToken / Celestia ---

$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D38") "2.82"
$ws.Range("E38").Value = "  -0.61%  "

$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D39") "3.21"
$ws.Range("E39").Value = "  -3.53%  "

$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue $ws.Range("D40") "18.33"
$ws.Range("E40").Value = "  -2.55%  "

# --- continue simple updates ---

$ws.Range("E41").Value = "  -3.33%  "

$ws.Range("E42").Value = "  +2.89%  "

Set-TextValue $ws.Range("D43") "22.98"
$ws.Range("E43").Value = "  -2.27%  "

Set-TextValue $ws.Range("D44") "120.17"
$ws.Range("E44").Value = "  -1.01%  "

Set-TextValue $ws.Range("D45") "2.17"
$ws.Range("E45").Value = "  -0.95%  "

# --- Row 46 / 47 swap (NEARProtocol <-> ApeXProtocol) with updated data ---

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D46") "2.48"
$ws.Range("E46").Value = "  -5.53%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D47") "3.47"
$ws.Range("E47").Value = "  -2.12%  "

# --- final simple updates ---

Set-TextValue $ws.Range("D48") "2.133.47"
$ws.Range("E48").Value = "  -2.88%  "

Set-TextValue $ws.Range("D49") "0.250"
$ws.Range("E49").Value = "  -4.98%  "

Set-TextValue $ws.Range("D50") "0.0349"
$ws.Range("E50").Value = "  +3.19%  "

Set-TextValue $ws.Range("D51") "0.933"
$ws.Range("E51").Value = "  -3.85%  "
